# Add a new "#ProductionPeriod" parameter block (rows 27-29) to the
# Vanern configuration-inputs sheet: ProdStartDay / ProdEndDay, both
# expressed in JulianDay units.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new labels in this order (#ProductionPeriod, ProdStartDay,
# ProdEndDay, then JulianDay) so the shared-string table is appended in
# the same sequence as the source workbook.
$ws.Range("A27").Value = "#ProductionPeriod"
$ws.Range("A28").Value = "ProdStartDay"
$ws.Range("A29").Value = "ProdEndDay"
$ws.Range("C28").Value = "JulianDay"

$ws.Range("B28").Value = 100
$ws.Range("C29").Value = "JulianDay"
$ws.Range("B29").Value = 300

# Match the author's scroll position / selection captured in the diff.
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("B30").Select()
